# Check for duplicate records during adding of records:
#  - read the existing contact table
#  - drop the Book_Name / Contact_Type columns
#  - de-duplicate records that share the same first_name, last_name,
#    address, city, state, zip and phone_number
#  - add a new "email" column with each contact's email address

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

# Read the full table (header row + data rows) using Value2, which (unlike
# the bare Value getter) reliably returns the underlying cell value here.
$table = @()
for ($r = 1; $r -le $rowCount; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $colCount; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $table += ,$rowVals
}

$headerRow = $table[0]

# Locate the key columns by header name so this does not depend on a fixed
# column layout.
$colIndex = @{}
for ($c = 0; $c -lt $headerRow.Count; $c++) {
    $colIndex[[string]$headerRow[$c]] = $c
}

$keyFields = @("first_name","last_name","address","city","state","zip","phone_number")

# Email address for each contact, looked up by first_name|last_name.
$emails = @{
    "Sue|Black"      = "sueblack@gmail.com"
    "Mike|Brown"     = "mikeybrown@gmail.com"
    "Liz|White"      = "lizwhite@yahoo.com"
    "Meg|Stephenson" = "megsteph@gmaco.com"
    "John|Gray"      = "johngray@gmail.com"
}

# Walk the data rows, checking for duplicate records (matching on the key
# fields) before adding each one to the de-duplicated list.
$seen = @{}
$uniqueRows = @()
for ($r = 1; $r -lt $table.Count; $r++) {
    $row = $table[$r]
    $keyParts = @()
    foreach ($field in $keyFields) {
        $keyParts += [string]$row[$colIndex[$field]]
    }
    $key = ($keyParts -join "|")
    if (-not $seen.ContainsKey($key)) {
        $seen[$key] = $true
        $uniqueRows += ,$row
    }
}

# Clear the sheet and rewrite it with the key fields plus a new email column
# (Book_Name / Contact_Type are dropped).
$ws.Cells.Clear()

$newHeaders = $keyFields + @("email")
for ($c = 0; $c -lt $newHeaders.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $newHeaders[$c]
}

for ($r = 0; $r -lt $uniqueRows.Count; $r++) {
    $row = $uniqueRows[$r]
    $excelRow = $r + 2
    $first = [string]$row[$colIndex["first_name"]]
    $last = [string]$row[$colIndex["last_name"]]
    for ($c = 0; $c -lt $keyFields.Count; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$colIndex[$keyFields[$c]]]
    }
    $ws.Cells.Item($excelRow, $keyFields.Count + 1).Value = $emails["$first|$last"]
}

$ws.Range("A1:J6").Select()
